$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.06754433333333333
$ws.Range("H2").Value2 = 0.202633
$ws.Range("I2").Value2 = 0.006855017925354449
$ws.Range("J2").Value2 = 0.006855017925354449
$ws.Range("M2").Value2 = 4.407279333333333
$ws.Range("N2").Value2 = 13.221838
$ws.Range("O2").Value2 = 0.2946616623342344
$ws.Range("P2").Value2 = 0.2946616623342344
$ws.Range("Q2").Value2 = 0.2976867443837778
$ws.Range("R2").Value2 = 2.679180699454
$ws.Range("S2").Value2 = 0.002019910977215916
$ws.Range("T2").Value2 = 0.002019910977215916
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.06754433333333333
$ws.Range("H3").Value2 = 0.202633
$ws.Range("I3").Value2 = 0.006855017925354449
$ws.Range("J3").Value2 = 0.006855017925354449
$ws.Range("O3").Value2 = 0.2393683991842171
$ws.Range("P3").Value2 = 0.2393683991842171
$ws.Range("Q3").Value2 = 0.2418258245644444
$ws.Range("R3").Value2 = 2.17643242108
$ws.Range("S3").Value2 = 0.001640874667171207
$ws.Range("T3").Value2 = 0.001640874667171207
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.06754433333333333
$ws.Range("H4").Value2 = 0.202633
$ws.Range("I4").Value2 = 0.006855017925354449
$ws.Range("J4").Value2 = 0.006855017925354449
$ws.Range("M4").Value2 = 3.580339
$ws.Range("N4").Value2 = 10.741017
$ws.Range("O4").Value2 = 0.2393741266819538
$ws.Range("P4").Value2 = 0.2393741266819538
$ws.Range("Q4").Value2 = 0.2418316108623333
$ws.Range("R4").Value2 = 2.176484497761
$ws.Range("S4").Value2 = 0.00164091392927086
$ws.Range("T4").Value2 = 0.00164091392927086
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.06754433333333333
$ws.Range("H5").Value2 = 0.202633
$ws.Range("I5").Value2 = 0.006855017925354449
$ws.Range("J5").Value2 = 0.006855017925354449
$ws.Range("M5").Value2 = 3.389212666666667
$ws.Range("N5").Value2 = 10.167638
$ws.Range("O5").Value2 = 0.2265958117995947
$ws.Range("P5").Value2 = 0.2265958117995947
$ws.Range("Q5").Value2 = 0.2289221100948889
$ws.Range("R5").Value2 = 2.060298990854
$ws.Range("S5").Value2 = 0.001553318351696465
$ws.Range("T5").Value2 = 0.001553318351696465
$ws.Range("I6").Value2 = 0.7774992501642265
$ws.Range("J6").Value2 = 0.7774992501642265
$ws.Range("M6").Value2 = 4.407279333333333
$ws.Range("N6").Value2 = 13.221838
$ws.Range("O6").Value2 = 0.2946616623342344
$ws.Range("P6").Value2 = 0.2946616623342344
$ws.Range("Q6").Value2 = 33.76376591024734
$ws.Range("R6").Value2 = 303.873893192226
$ws.Range("S6").Value2 = 0.2290992215170117
$ws.Range("T6").Value2 = 0.2290992215170117
$ws.Range("I7").Value2 = 0.7774992501642265
$ws.Range("J7").Value2 = 0.7774992501642265
$ws.Range("O7").Value2 = 0.2393683991842171
$ws.Range("P7").Value2 = 0.2393683991842171
$ws.Range("S7").Value2 = 0.18610875087874
$ws.Range("T7").Value2 = 0.18610875087874
$ws.Range("I8").Value2 = 0.7774992501642265
$ws.Range("J8").Value2 = 0.7774992501642265
$ws.Range("M8").Value2 = 3.580339
$ws.Range("N8").Value2 = 10.741017
$ws.Range("O8").Value2 = 0.2393741266819538
$ws.Range("P8").Value2 = 0.2393741266819538
$ws.Range("Q8").Value2 = 27.428651268151
$ws.Range("R8").Value2 = 246.857861413359
$ws.Range("S8").Value2 = 0.1861132040039356
$ws.Range("T8").Value2 = 0.1861132040039357
$ws.Range("I9").Value2 = 0.7774992501642265
$ws.Range("J9").Value2 = 0.7774992501642265
$ws.Range("M9").Value2 = 3.389212666666667
$ws.Range("N9").Value2 = 10.167638
$ws.Range("O9").Value2 = 0.2265958117995947
$ws.Range("P9").Value2 = 0.2265958117995947
$ws.Range("Q9").Value2 = 25.96444982098067
$ws.Range("R9").Value2 = 233.680048388826
$ws.Range("S9").Value2 = 0.1761780737645391
$ws.Range("T9").Value2 = 0.1761780737645391
$ws.Range("G10").Value2 = 1.941983333333333
$ws.Range("H10").Value2 = 5.825949999999999
$ws.Range("I10").Value2 = 0.1970902650714284
$ws.Range("J10").Value2 = 0.1970902650714283
$ws.Range("M10").Value2 = 4.407279333333333
$ws.Range("N10").Value2 = 13.221838
$ws.Range("O10").Value2 = 0.2946616623342344
$ws.Range("P10").Value2 = 0.2946616623342344
$ws.Range("Q10").Value2 = 8.558863010677776
$ws.Range("R10").Value2 = 77.02976709609999
$ws.Range("S10").Value2 = 0.05807494513584197
$ws.Range("T10").Value2 = 0.05807494513584197
$ws.Range("G11").Value2 = 1.941983333333333
$ws.Range("H11").Value2 = 5.825949999999999
$ws.Range("I11").Value2 = 0.1970902650714284
$ws.Range("J11").Value2 = 0.1970902650714283
$ws.Range("O11").Value2 = 0.2393683991842171
$ws.Range("P11").Value2 = 0.2393683991842171
$ws.Range("Q11").Value2 = 6.952792302444442
$ws.Range("R11").Value2 = 62.57513072199998
$ws.Range("S11").Value2 = 0.04717718124494082
$ws.Range("T11").Value2 = 0.04717718124494082
$ws.Range("G12").Value2 = 1.941983333333333
$ws.Range("H12").Value2 = 5.825949999999999
$ws.Range("I12").Value2 = 0.1970902650714284
$ws.Range("J12").Value2 = 0.1970902650714283
$ws.Range("M12").Value2 = 3.580339
$ws.Range("N12").Value2 = 10.741017
$ws.Range("O12").Value2 = 0.2393741266819538
$ws.Range("P12").Value2 = 0.2393741266819538
$ws.Range("Q12").Value2 = 6.952958665683332
$ws.Range("R12").Value2 = 62.57662799114998
$ws.Range("S12").Value2 = 0.04717831007898795
$ws.Range("T12").Value2 = 0.04717831007898795
$ws.Range("G13").Value2 = 1.941983333333333
$ws.Range("H13").Value2 = 5.825949999999999
$ws.Range("I13").Value2 = 0.1970902650714284
$ws.Range("J13").Value2 = 0.1970902650714283
$ws.Range("M13").Value2 = 3.389212666666667
$ws.Range("N13").Value2 = 10.167638
$ws.Range("O13").Value2 = 0.2265958117995947
$ws.Range("P13").Value2 = 0.2265958117995947
$ws.Range("Q13").Value2 = 6.581794511788888
$ws.Range("R13").Value2 = 59.23615060609999
$ws.Range("S13").Value2 = 0.04465982861165763
$ws.Range("T13").Value2 = 0.04465982861165762
$ws.Range("E14").Value2 = 2
$ws.Range("F14").Value2 = 0.6666666666666666
$ws.Range("G14").Value2 = 0.182832
$ws.Range("H14").Value2 = 0.548496
$ws.Range("I14").Value2 = 0.01855546683899075
$ws.Range("J14").Value2 = 0.01855546683899075
$ws.Range("M14").Value2 = 4.407279333333333
$ws.Range("N14").Value2 = 13.221838
$ws.Range("O14").Value2 = 0.2946616623342344
$ws.Range("P14").Value2 = 0.2946616623342344
$ws.Range("Q14").Value2 = 0.805791695072
$ws.Range("R14").Value2 = 7.252125255648
$ws.Range("S14").Value2 = 0.005467584704164777
$ws.Range("T14").Value2 = 0.005467584704164777
$ws.Range("E15").Value2 = 2
$ws.Range("F15").Value2 = 0.6666666666666666
$ws.Range("G15").Value2 = 0.182832
$ws.Range("H15").Value2 = 0.548496
$ws.Range("I15").Value2 = 0.01855546683899075
$ws.Range("J15").Value2 = 0.01855546683899075
$ws.Range("O15").Value2 = 0.2393683991842171
$ws.Range("P15").Value2 = 0.2393683991842171
$ws.Range("Q15").Value2 = 0.6545848774399999
$ws.Range("R15").Value2 = 5.891263896959999
$ws.Range("S15").Value2 = 0.004441592393365041
$ws.Range("T15").Value2 = 0.004441592393365042
$ws.Range("E16").Value2 = 2
$ws.Range("F16").Value2 = 0.6666666666666666
$ws.Range("G16").Value2 = 0.182832
$ws.Range("H16").Value2 = 0.548496
$ws.Range("I16").Value2 = 0.01855546683899075
$ws.Range("J16").Value2 = 0.01855546683899075
$ws.Range("M16").Value2 = 3.580339
$ws.Range("N16").Value2 = 10.741017
$ws.Range("O16").Value2 = 0.2393741266819538
$ws.Range("P16").Value2 = 0.2393741266819538
$ws.Range("Q16").Value2 = 0.6546005400479999
$ws.Range("R16").Value2 = 5.891404860431999
$ws.Range("S16").Value2 = 0.004441698669759366
$ws.Range("T16").Value2 = 0.004441698669759366
$ws.Range("E17").Value2 = 2
$ws.Range("F17").Value2 = 0.6666666666666666
$ws.Range("G17").Value2 = 0.182832
$ws.Range("H17").Value2 = 0.548496
$ws.Range("I17").Value2 = 0.01855546683899075
$ws.Range("J17").Value2 = 0.01855546683899075
$ws.Range("M17").Value2 = 3.389212666666667
$ws.Range("N17").Value2 = 10.167638
$ws.Range("O17").Value2 = 0.2265958117995947
$ws.Range("P17").Value2 = 0.2265958117995947
$ws.Range("Q17").Value2 = 0.6196565302720001
$ws.Range("R17").Value2 = 5.576908772448
$ws.Range("S17").Value2 = 0.00420459107170157
$ws.Range("T17").Value2 = 0.00420459107170157
